# "CAJAS PARA ENRROLLADOR" price list refresh:
#  - bump the quoted date in A1 to the new list date
#  - update the two enrollador box prices

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A1 holds the price-list date as an Excel serial number (formatted as a date).
# 45406 -> 2024-04-24, 45436 -> 2024-05-24
$ws.Range("A1").Value = 45436

# CAJA p/ ENROLLADOR CHICA
$ws.Range("D30").Value = 799

# CAJA p/ ENROLLADOR GRANDE
$ws.Range("D31").Value = 967
